$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number, new date text (dashes instead of slashes),
# and the new D/E/F/G/H values.
$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 11; Date = "25-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 12; Date = "29-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 13; Date = "01-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 14; Date = "05-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 15; Date = "08-09-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 16; Date = "12-09-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 21; Date = "29-09-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 }
)

# Rows whose DD-MM-YYYY text is ambiguous (day <= 12) and would otherwise
# be auto-parsed into a date serial number by plain .Value assignment.
# Only these need a leading apostrophe to stay as literal text, keeping the
# rest of the cells (and their style index) untouched, matching the source
# edit exactly.
$ambiguousRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($r in $rows) {
    $cellA = $ws.Cells.Item($r.Row, 1)
    if ($ambiguousRows -contains $r.Row) {
        $cellA.Value = "'" + $r.Date
    } else {
        $cellA.Value = $r.Date
    }

    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
